$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.194193959236145
$ws.Range("B1").Value = 2.381388187408447
$ws.Range("C1").Value = 4.318923473358154
$ws.Range("D1").Value = 2.76301646232605
$ws.Range("E1").Value = 1.115826845169067
